$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = "0003"
$ws.Range("B4").Value = "Fling Match Lighter Metal Outdoor Camping"
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 149
$ws.Range("E4").Value = 1
